$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header C1 from mean_pCO2 to median_pCO2
$ws.Range("C1").Value = "median_pCO2"

$data = @{
    2 = @{ B = 354.312659315757; D = 1817.7522476787 }
    3 = @{ B = 374.984834279778; D = 1851.59427516284 }
    4 = @{ B = 353.892957746479; D = 1630.94332552693 }
    5 = @{ B = 216.137688666002; D = 995.22224953073 }
    6 = @{ B = 135.378030604928; D = 731.123832624966 }
    7 = @{ B = 206.093685390115; D = 1067.85681923597 }
    8 = @{ B = 529.340625566246; D = 2782.32020958726 }
    9 = @{ B = 156.569767658106; D = 604.498651462128 }
    10 = @{ B = 156.845958742782; D = 598.975443403997 }
    11 = @{ B = 229.669238439386; D = 851.646775882724 }
    12 = @{ B = 260.363550388584; D = 1020.36513260567 }
    13 = @{ B = 449.507402408847; D = 1692.68435682913 }
    14 = @{ B = 280.030289857542; D = 1078.96427461686 }
    15 = @{ B = 557.211232419769; D = 2068.78947817386 }
    16 = @{ B = 249.89123999692; D = 1058.60754177337 }
    17 = @{ B = 533.363431057008; D = 1999.49913475042 }
    18 = @{ B = 538.942825535404; D = 2072.19980004882 }
    19 = @{ B = 454.347461223283; D = 1720.34194977037 }
    20 = @{ B = 469.525406671454; D = 1854.2900208017 }
    21 = @{ B = 487.296760859198; D = 1895.25200954586 }
    22 = @{ B = 216.81281391383; D = 791.936239372724 }
    23 = @{ B = 517.766977277684; D = 2191.995183028 }
    24 = @{ B = 715.061720362547; D = 2844.44689668828 }
    25 = @{ B = 667.206057120518; D = 2871.28722931183 }
    26 = @{ B = 289.766781433583; D = 1207.68538137971 }
    27 = @{ B = 593.088632108031; D = 2719.41827785334 }
    28 = @{ B = 412.096366237632; D = 1597.50328149219 }
    29 = @{ B = 191.279079447917; D = 792.971626044365 }
    30 = @{ B = 206.221650425231; D = 1005.01159361935 }
    31 = @{ B = 203.144630664703; D = 847.382655260253 }
    32 = @{ B = 645.892587233709; D = 2338.47748702646 }
    33 = @{ B = 841.705442971699; D = 2970.18441440874 }
    34 = @{ B = 660.30475218755; D = 2251.361181084 }
    35 = @{ B = 550.505207716531; D = 1934.18646323408 }
    36 = @{ B = 1124.85721896266; D = 4006.65169540156 }
    37 = @{ B = 1119.98811416649; D = 3981.65866647652 }
    38 = @{ B = 978.755625748932; D = 3443.67624577205 }
    39 = @{ B = 943.291487066612; D = 3347.36304645877 }
    40 = @{ B = 665.386291658186; D = 2338.17295116595 }
    41 = @{ B = 610.166617640876; D = 2205.05439890485 }
    42 = @{ B = 141.581050401698; D = 552.858823354043 }
    43 = @{ B = 144.189712427675; D = 561.964627310409 }
    44 = @{ B = 600.540849909865; D = 2037.62729743552 }
    45 = @{ B = 569.196264007643; D = 2011.0768337519 }
    46 = @{ B = 221.115902897825; D = 869.794315577981 }
    47 = @{ B = 569.291456582512; D = 2591.78985041937 }
    48 = @{ B = 566.526894936938; D = 2579.80000055119 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 4).Value = $vals.D
}
